$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 957.7143
$ws.Range("I6").Value = 140.8
$ws.Range("K6").Value = 422.4
$ws.Range("M6").Value = -310.4
$ws.Range("H19").Value = 408.47827
$ws.Range("I19").Value = 451.18182
$ws.Range("J19").Value = 369.33334
$ws.Range("K19").Value = 451.18182
$ws.Range("L19").Value = 369.33334
$ws.Range("M19").Value = -276.18182
$ws.Range("N19").Value = -719.33334
$ws.Range("H40").Value = 4401.75
$ws.Range("I40").Value = 3498.7334
$ws.Range("K40").Value = 3498.7334
$ws.Range("M40").Value = -3323.7334
$ws.Range("H53").Value = 1748.875
$ws.Range("I53").Value = 825.5
$ws.Range("K53").Value = 825.5
$ws.Range("M53").Value = -188.5
$ws.Range("H74").Value = 12973.167
$ws.Range("I74").Value = 8669.429
$ws.Range("J74").Value = 18998.4
$ws.Range("K74").Value = 8669.429
$ws.Range("L74").Value = 18998.4
$ws.Range("M74").Value = -7733.429
$ws.Range("N74").Value = -20870.4
$ws.Range("H77").Value = 12973.167
$ws.Range("I77").Value = 8669.429
$ws.Range("J77").Value = 18998.4
$ws.Range("K77").Value = 43347.145
$ws.Range("L77").Value = 94992
$ws.Range("M77").Value = -38667.145
$ws.Range("N77").Value = -104352
$ws.Range("H118").Value = 1595.1428
$ws.Range("I118").Value = 433.2
$ws.Range("J118").Value = 4500
$ws.Range("K118").Value = 1299.6
$ws.Range("L118").Value = 13500
$ws.Range("M118").Value = 357.4000000000001
$ws.Range("N118").Value = -16814
$ws.Range("H125").Value = 5000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("H138").Value = 3179.8
$ws.Range("I138").Value = 2890.6667
$ws.Range("K138").Value = 8672.000100000001
$ws.Range("M138").Value = -3532.000100000001
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2520.919
$ws.Range("I32").Value = 2520.919
$ws.Range("K32").Value = 2520.919
$ws.Range("M32").Value = -2233.919
$ws.Range("H61").Value = 6050
$ws.Range("I61").Value = 4600
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 4600
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -4388
$ws.Range("N61").Value = -7924
$ws.Range("H63").Value = 5118.4
$ws.Range("I63").Value = 2794
$ws.Range("K63").Value = 2794
$ws.Range("M63").Value = -2108
$ws.Range("H66").Value = 5118.4
$ws.Range("I66").Value = 2794
$ws.Range("K66").Value = 13970
$ws.Range("M66").Value = -10538
$ws.Range("H80").Value = 50082.5
$ws.Range("J80").Value = 50082.5
$ws.Range("L80").Value = 50082.5
$ws.Range("N80").Value = -52078.5
$ws.Range("H83").Value = 50082.5
$ws.Range("J83").Value = 50082.5
$ws.Range("L83").Value = 150247.5
$ws.Range("N83").Value = -160231.5
$ws.Range("H136").Value = 6050
$ws.Range("I136").Value = 4600
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 13800
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -11250
$ws.Range("N136").Value = -27600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 44051.855
$ws.Range("I82").Value = 19678.4
$ws.Range("J82").Value = 104985.5
$ws.Range("K82").Value = 19678.4
$ws.Range("L82").Value = 104985.5
$ws.Range("M82").Value = -19295.4
$ws.Range("N82").Value = -105751.5
$ws.Range("H85").Value = 44051.855
$ws.Range("I85").Value = 19678.4
$ws.Range("J85").Value = 104985.5
$ws.Range("K85").Value = 19678.4
$ws.Range("L85").Value = 104985.5
$ws.Range("M85").Value = -18352.4
$ws.Range("N85").Value = -107637.5
$ws.Range("H86").Value = 6060.9165
$ws.Range("I86").Value = 5273.1
$ws.Range("K86").Value = 5273.1
$ws.Range("M86").Value = -4150.1
$ws.Range("H89").Value = 6060.9165
$ws.Range("I89").Value = 5273.1
$ws.Range("K89").Value = 26365.5
$ws.Range("M89").Value = -20749.5
$ws.Range("H105").Value = 2077.6
$ws.Range("I105").Value = 2077.6
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2077.6
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -330.5999999999999
$ws.Range("H134").Value = 6419.1113
$ws.Range("I134").Value = 5971.5
$ws.Range("K134").Value = 17914.5
$ws.Range("M134").Value = -15379.5
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 550
$ws.Range("I14").Value = 550
$ws.Range("K14").Value = 550
$ws.Range("M14").Value = -380
$ws.Range("H16").Value = 1885.3334
$ws.Range("I16").Value = 1754.875
$ws.Range("J16").Value = 2929
$ws.Range("K16").Value = 1754.875
$ws.Range("L16").Value = 2929
$ws.Range("M16").Value = -1467.875
$ws.Range("N16").Value = -3503
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
$ws.Range("H105").Value = 2769.1667
$ws.Range("I105").Value = 2434.2
$ws.Range("K105").Value = 2434.2
$ws.Range("M105").Value = -687.1999999999998
$ws.Range("H113").Value = 1885.3334
$ws.Range("I113").Value = 1754.875
$ws.Range("J113").Value = 2929
$ws.Range("K113").Value = 1754.875
$ws.Range("L113").Value = 2929
$ws.Range("M113").Value = 415.125
$ws.Range("N113").Value = -7269
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6764.9414
$ws.Range("J39").Value = 8691.308000000001
$ws.Range("L39").Value = 26073.924
$ws.Range("N39").Value = -26661.924
$ws.Range("H80").Value = 3888.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3888.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 11665.5
$ws.Range("N80").Value = -13537.5
$ws.Range("H83").Value = 3888.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3888.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 34996.5
$ws.Range("N83").Value = -44356.5
$ws.Range("H132").Value = 2099.8572
$ws.Range("I132").Value = 2033.1666
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 18298.4994
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -15768.4994
$ws.Range("N132").Value = -27560
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 126.125
$ws.Range("J2").Value = 175
$ws.Range("L2").Value = 175
$ws.Range("N2").Value = -401
$ws.Range("H80").Value = 3219.2
$ws.Range("I80").Value = 3024.5
$ws.Range("K80").Value = 3024.5
$ws.Range("M80").Value = -2026.5
$ws.Range("H83").Value = 3219.2
$ws.Range("I83").Value = 3024.5
$ws.Range("K83").Value = 15122.5
$ws.Range("M83").Value = -10130.5
$ws.Range("H86").Value = 35000
$ws.Range("J86").Value = 35000
$ws.Range("L86").Value = 35000
$ws.Range("N86").Value = -37372
$ws.Range("H89").Value = 35000
$ws.Range("J89").Value = 35000
$ws.Range("L89").Value = 105000
$ws.Range("N89").Value = -116856
$ws.Range("H113").Value = 9774.4
$ws.Range("I113").Value = 9249.666999999999
$ws.Range("J113").Value = 9999.286
$ws.Range("K113").Value = 9249.666999999999
$ws.Range("L113").Value = 9999.286
$ws.Range("M113").Value = -7079.666999999999
$ws.Range("N113").Value = -14339.286
$ws.Range("H126").Value = 3943.8667
$ws.Range("I126").Value = 3704.4614
$ws.Range("K126").Value = 11113.3842
$ws.Range("M126").Value = -8643.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6275.421
$ws.Range("I40").Value = 4955.643
$ws.Range("J40").Value = 9970.799999999999
$ws.Range("K40").Value = 4955.643
$ws.Range("L40").Value = 9970.799999999999
$ws.Range("M40").Value = -4819.643
$ws.Range("N40").Value = -10242.8
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 207
$ws.Range("J18").Value = 207
$ws.Range("L18").Value = 207
$ws.Range("N18").Value = -553
$ws.Range("H57").Value = 100000
$ws.Range("J57").Value = 100000
$ws.Range("L57").Value = 100000
$ws.Range("N57").Value = -101508
$ws.Range("H100").Value = 489.53845
$ws.Range("I100").Value = 435.66666
$ws.Range("K100").Value = 871.33332
$ws.Range("M100").Value = -330.33332
$ws.Range("H126").Value = 3910.611
$ws.Range("I126").Value = 1866
$ws.Range("K126").Value = 5598
$ws.Range("M126").Value = -3128
